$wb = $excel.ActiveWorkbook

# Sheet "展览" (index/rId1) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8865
$ws1.Range("F3").Value = 8295
$ws1.Range("F4").Value = 145
$ws1.Range("F13").Value = 208
$ws1.Range("F14").Value = 5370
$ws1.Range("F18").Value = 14
$ws1.Range("F21").Value = 157
$ws1.Range("F22").Value = 177
$ws1.Range("F23").Value = 11

# Sheet "全部类型" (index/rId4) - same rows, F22 differs from sheet "展览"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8865
$ws4.Range("F3").Value = 8295
$ws4.Range("F4").Value = 145
$ws4.Range("F13").Value = 208
$ws4.Range("F14").Value = 5370
$ws4.Range("F18").Value = 14
$ws4.Range("F21").Value = 157
$ws4.Range("F22").Value = 178
$ws4.Range("F23").Value = 11
